# "sample num 4 is cleaned up"
#
# On the Aircraft_scheduling sheet, several rows in column E held a
# formula of the form "=C<n>+<offset>" whose result crossed midnight
# (i.e. landed above 1.0, a full extra day), or otherwise carried long
# binary-floating-point noise. Those cells are cleaned up by replacing
# the live formula with its already-reduced (mod 1 day) literal time
# value, so the sheet stores a plain number instead of a formula.
# Column F (elapsed-time, a shared MOD(E-C,1) formula) naturally picks
# up the freshly reduced values once E changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aircraft_scheduling")

$ws.Range("E4").Value  = 0.1711111111111111
$ws.Range("E5").Value  = 0.13930555555555554
$ws.Range("E6").Value  = 0.8666898148148148
$ws.Range("E9").Value  = 0.1625
$ws.Range("E11").Value = 0.3154398148148148
$ws.Range("E12").Value = 0.31916666666666665
$ws.Range("E13").Value = 0.43277777777777776
$ws.Range("E14").Value = 0.6740277777777778
$ws.Range("E15").Value = 0.47907407407407404
$ws.Range("E16").Value = 0.6039699074074074
$ws.Range("E17").Value = 0.5876388888888889
$ws.Range("E18").Value = 0.7354166666666666

# Cursor/selection ends up on E19 on this sheet.
$ws.Activate() | Out-Null
$ws.Range("E19").Select() | Out-Null
